$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67634766666667
$ws.Range("H2").Value = 83.029043
$ws.Range("I2").Value = 0.005965811625935536
$ws.Range("J2").Value = 0.005965811625935536
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 53.10831882110189
$ws.Range("R2").Value = 477.974869389917
$ws.Range("S2").Value = 0.00003892106936752136
$ws.Range("T2").Value = 0.00003892106936752136
$ws.Range("G3").Value = 27.67634766666667
$ws.Range("H3").Value = 83.029043
$ws.Range("I3").Value = 0.005965811625935536
$ws.Range("J3").Value = 0.005965811625935536
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 5017.400546472054
$ws.Range("R3").Value = 45156.60491824849
$ws.Range("S3").Value = 0.003677062257830038
$ws.Range("T3").Value = 0.003677062257830038
$ws.Range("G4").Value = 27.67634766666667
$ws.Range("H4").Value = 83.029043
$ws.Range("I4").Value = 0.005965811625935536
$ws.Range("J4").Value = 0.005965811625935536
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 803.933478404748
$ws.Range("R4").Value = 7235.401305642731
$ws.Range("S4").Value = 0.0005891723062307007
$ws.Range("T4").Value = 0.0005891723062307008
$ws.Range("G5").Value = 27.67634766666667
$ws.Range("H5").Value = 83.029043
$ws.Range("I5").Value = 0.005965811625935536
$ws.Range("J5").Value = 0.005965811625935536
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 2265.987274641688
$ws.Range("R5").Value = 20393.88547177519
$ws.Range("S5").Value = 0.001660655992507276
$ws.Range("T5").Value = 0.001660655992507276
$ws.Range("I6").Value = 0.009118181457976757
$ws.Range("J6").Value = 0.009118181457976757
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 81.17106578318266
$ws.Range("R6").Value = 730.539592048644
$ws.Range("S6").Value = 0.00005948719055907301
$ws.Range("T6").Value = 0.00005948719055907301
$ws.Range("I7").Value = 0.009118181457976757
$ws.Range("J7").Value = 0.009118181457976757
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.005620043508147854
$ws.Range("T7").Value = 0.005620043508147854
$ws.Range("I8").Value = 0.009118181457976757
$ws.Range("J8").Value = 0.009118181457976757
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 1228.736640689255
$ws.Range("R8").Value = 11058.62976620329
$ws.Range("S8").Value = 0.0009004944063053171
$ws.Range("T8").Value = 0.0009004944063053172
$ws.Range("I9").Value = 0.009118181457976757
$ws.Range("J9").Value = 0.009118181457976757
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 3463.34823275769
$ws.Range("R9").Value = 31170.13409481922
$ws.Range("S9").Value = 0.002538156352964513
$ws.Range("T9").Value = 0.002538156352964513
$ws.Range("G10").Value = 29.593002
$ws.Range("H10").Value = 88.779006
$ws.Range("I10").Value = 0.006378958578792732
$ws.Range("J10").Value = 0.006378958578792732
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 56.786198960146
$ws.Range("R10").Value = 511.0757906413139
$ws.Range("S10").Value = 0.00004161644800489384
$ws.Range("T10").Value = 0.00004161644800489384
$ws.Range("G11").Value = 29.593002
$ws.Range("H11").Value = 88.779006
$ws.Range("I11").Value = 0.006378958578792732
$ws.Range("J11").Value = 0.006378958578792732
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 5364.867727304117
$ws.Range("R11").Value = 48283.80954573704
$ws.Range("S11").Value = 0.003931707754963121
$ws.Range("T11").Value = 0.003931707754963121
$ws.Range("G12").Value = 29.593002
$ws.Range("H12").Value = 88.779006
$ws.Range("I12").Value = 0.006378958578792732
$ws.Range("J12").Value = 0.006378958578792732
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 859.607825455678
$ws.Range("R12").Value = 7736.470429101101
$ws.Range("S12").Value = 0.00062997392020873
$ws.Range("T12").Value = 0.0006299739202087301
$ws.Range("G13").Value = 29.593002
$ws.Range("H13").Value = 88.779006
$ws.Range("I13").Value = 0.006378958578792732
$ws.Range("J13").Value = 0.006378958578792732
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 2422.912400078344
$ws.Range("R13").Value = 21806.21160070509
$ws.Range("S13").Value = 0.001775660455615987
$ws.Range("T13").Value = 0.001775660455615987
$ws.Range("G14").Value = 4539.588785666667
$ws.Range("H14").Value = 13618.766357
$ws.Range("I14").Value = 0.9785370483372949
$ws.Range("J14").Value = 0.978537048337295
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 8711.045671544742
$ws.Range("R14").Value = 78399.41104390268
$ws.Range("S14").Value = 0.006383994454577335
$ws.Range("T14").Value = 0.006383994454577336
$ws.Range("G15").Value = 4539.588785666667
$ws.Range("H15").Value = 13618.766357
$ws.Range("I15").Value = 0.9785370483372949
$ws.Range("J15").Value = 0.978537048337295
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 822974.7482683503
$ws.Range("R15").Value = 7406772.734415152
$ws.Range("S15").Value = 0.6031269295676475
$ws.Range("T15").Value = 0.6031269295676476
$ws.Range("G16").Value = 4539.588785666667
$ws.Range("H16").Value = 13618.766357
$ws.Range("I16").Value = 0.9785370483372949
$ws.Range("J16").Value = 0.978537048337295
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 131864.4875741199
$ws.Range("R16").Value = 1186780.388167079
$ws.Range("S16").Value = 0.09663847363109759
$ws.Range("T16").Value = 0.09663847363109762
$ws.Range("G17").Value = 4539.588785666667
$ws.Range("H17").Value = 13618.766357
$ws.Range("I17").Value = 0.9785370483372949
$ws.Range("J17").Value = 0.978537048337295
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 371676.5862432057
$ws.Range("R17").Value = 3345089.276188851
$ws.Range("S17").Value = 0.2723876506839725
$ws.Range("T17").Value = 0.2723876506839725
